# Update STATUS / COMMENT values in John_Quest.xlsx
# - Sheet1: E2 OK -> NO ISSUE, F2 comment update
#           E3 ERROR -> ISSUE, F3 comment update
#           E5 OK -> NO ISSUE
# - Sheet2: E2 OK -> NO ISSUE
#           E4 PENDING -> BLOCKED

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("E2").Value = "NO ISSUE"
$ws1.Range("F2").Value = "Translation looks good"

$ws1.Range("E3").Value = "ISSUE"
$ws1.Range("F3").Value = "Typo: should be singular"

$ws1.Range("E5").Value = "NO ISSUE"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("E2").Value = "NO ISSUE"

$ws2.Range("E4").Value = "BLOCKED"
